$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 152.66667
$ws.Range("I2").Value = 119
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 119
$ws.Range("L2").Value = 220
$ws.Range("M2").Value = -6
$ws.Range("N2").Value = -446

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 613.24
$ws.Range("J19").Value = 592.5333000000001
$ws.Range("L19").Value = 592.5333000000001
$ws.Range("N19").Value = -942.5333000000001

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 109139.5
$ws.Range("I116").Value = 134861.88
$ws.Range("J116").Value = 6250
$ws.Range("K116").Value = 134861.88
$ws.Range("L116").Value = 6250
$ws.Range("M116").Value = -131419.88
$ws.Range("N116").Value = -13134

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 35394
$ws.Range("J140").Value = 35394
$ws.Range("L140").Value = 35394
$ws.Range("N140").Value = -45754

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6240.33
$ws.Range("I32").Value = 3971
$ws.Range("J32").Value = 18154.312
$ws.Range("K32").Value = 3971
$ws.Range("L32").Value = 18154.312
$ws.Range("M32").Value = -3684
$ws.Range("N32").Value = -18728.312

# ARM row 75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 31420
$ws.Range("J75").Value = 31420
$ws.Range("L75").Value = 31420
$ws.Range("N75").Value = -33168

# ARM row 78
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 31420
$ws.Range("J78").Value = 31420
$ws.Range("L78").Value = 94260
$ws.Range("N78").Value = -102996

# ARM row 82
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 31453.334
$ws.Range("J82").Value = 31453.334
$ws.Range("L82").Value = 31453.334
$ws.Range("N82").Value = -32175.334

# ARM row 85
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 31453.334
$ws.Range("J85").Value = 31453.334
$ws.Range("L85").Value = 31453.334
$ws.Range("N85").Value = -33949.334

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3975
$ws.Range("I88").Value = 4560
$ws.Range("K88").Value = 4560
$ws.Range("M88").Value = -4154

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3975
$ws.Range("I91").Value = 4560
$ws.Range("K91").Value = 4560
$ws.Range("M91").Value = -3156

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2671.1155
$ws.Range("I122").Value = 2815.3157
$ws.Range("K122").Value = 8445.947100000001
$ws.Range("M122").Value = -5995.947100000001

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7249.722
$ws.Range("I86").Value = 6712.8
$ws.Range("J86").Value = 7920.875
$ws.Range("K86").Value = 6712.8
$ws.Range("L86").Value = 7920.875
$ws.Range("M86").Value = -5589.8
$ws.Range("N86").Value = -10166.875

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7249.722
$ws.Range("I89").Value = 6712.8
$ws.Range("J89").Value = 7920.875
$ws.Range("K89").Value = 33564
$ws.Range("L89").Value = 39604.375
$ws.Range("M89").Value = -27948
$ws.Range("N89").Value = -50836.375

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1132.2858
$ws.Range("I107").Value = 1200.1666
$ws.Range("K107").Value = 1200.1666
$ws.Range("M107").Value = 719.8334

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 23683.334
$ws.Range("I134").Value = 30911.766
$ws.Range("J134").Value = 6128.5713
$ws.Range("K134").Value = 92735.298
$ws.Range("L134").Value = 18385.7139
$ws.Range("M134").Value = -90200.298
$ws.Range("N134").Value = -23455.7139

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 53693.21
$ws.Range("I16").Value = 100723.1
$ws.Range("J16").Value = 1437.7778
$ws.Range("K16").Value = 100723.1
$ws.Range("L16").Value = 1437.7778
$ws.Range("M16").Value = -100436.1
$ws.Range("N16").Value = -2011.7778

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3730.875
$ws.Range("I62").Value = 3504.1667
$ws.Range("J62").Value = 4411
$ws.Range("K62").Value = 3504.1667
$ws.Range("L62").Value = 4411
$ws.Range("M62").Value = -2880.1667
$ws.Range("N62").Value = -5659

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3730.875
$ws.Range("I65").Value = 3504.1667
$ws.Range("J65").Value = 4411
$ws.Range("K65").Value = 17520.8335
$ws.Range("L65").Value = 22055
$ws.Range("M65").Value = -14400.8335
$ws.Range("N65").Value = -28295

# CRP row 87
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 103388.336
$ws.Range("J87").Value = 103388.336
$ws.Range("L87").Value = 103388.336
$ws.Range("N87").Value = -105760.336

# CRP row 90
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 103388.336
$ws.Range("J90").Value = 103388.336
$ws.Range("L90").Value = 310165.008
$ws.Range("N90").Value = -322021.008

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 621.8033
$ws.Range("I105").Value = 611.85187
$ws.Range("J105").Value = 698.5714
$ws.Range("K105").Value = 611.85187
$ws.Range("L105").Value = 698.5714
$ws.Range("M105").Value = 1135.14813
$ws.Range("N105").Value = -4192.5714

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 359.04544
$ws.Range("I107").Value = 267.8125
$ws.Range("K107").Value = 267.8125
$ws.Range("M107").Value = 1652.1875

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 53693.21
$ws.Range("I113").Value = 100723.1
$ws.Range("J113").Value = 1437.7778
$ws.Range("K113").Value = 100723.1
$ws.Range("L113").Value = 1437.7778
$ws.Range("M113").Value = -98553.10000000001
$ws.Range("N113").Value = -5777.7778

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1747.4348
$ws.Range("I134").Value = 1284.8276
$ws.Range("J134").Value = 2536.5881
$ws.Range("K134").Value = 3854.4828
$ws.Range("L134").Value = 7609.7643
$ws.Range("M134").Value = -1319.4828
$ws.Range("N134").Value = -12679.7643

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 125.36364
$ws.Range("I8").Value = 125.36364
$ws.Range("K8").Value = 376.09092
$ws.Range("M8").Value = -237.09092

# CUL row 41
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1400
$ws.Range("J41").Value = 1950
$ws.Range("L41").Value = 5850
$ws.Range("N41").Value = -6526

# CUL row 60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 314.85715
$ws.Range("I60").Value = 110.888885
$ws.Range("J60").Value = 682
$ws.Range("K60").Value = 332.666655
$ws.Range("L60").Value = 2046
$ws.Range("M60").Value = -81.66665499999999
$ws.Range("N60").Value = -2548

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 5988.778
$ws.Range("I98").Value = 298.75
$ws.Range("J98").Value = 10540.8
$ws.Range("K98").Value = 896.25
$ws.Range("L98").Value = 31622.4
$ws.Range("M98").Value = 601.75
$ws.Range("N98").Value = -34618.39999999999

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 875.42426
$ws.Range("I131").Value = 693.6667
$ws.Range("J131").Value = 924.35895
$ws.Range("K131").Value = 2081.0001
$ws.Range("L131").Value = 2773.07685
$ws.Range("M131").Value = 2958.9999
$ws.Range("N131").Value = -12853.07685

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2668.2222
$ws.Range("I137").Value = 1630.3846
$ws.Range("J137").Value = 5366.6
$ws.Range("K137").Value = 4891.1538
$ws.Range("L137").Value = 16099.8
$ws.Range("M137").Value = 208.8462
$ws.Range("N137").Value = -26299.8

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11725.0625
$ws.Range("I70").Value = 3966.8333
$ws.Range("K70").Value = 3966.8333
$ws.Range("M70").Value = -3696.8333

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11725.0625
$ws.Range("I73").Value = 3966.8333
$ws.Range("K73").Value = 3966.8333
$ws.Range("M73").Value = -3030.8333

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9344.727999999999
$ws.Range("I107").Value = 12676.125
$ws.Range("J107").Value = 461
$ws.Range("K107").Value = 12676.125
$ws.Range("L107").Value = 461
$ws.Range("M107").Value = -10756.125
$ws.Range("N107").Value = -4301

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3475.1177
$ws.Range("I126").Value = 3837.375
$ws.Range("J126").Value = 3153.111
$ws.Range("K126").Value = 11512.125
$ws.Range("L126").Value = 9459.332999999999
$ws.Range("M126").Value = -9042.125
$ws.Range("N126").Value = -14399.333

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2893.3684
$ws.Range("I132").Value = 2696.889
$ws.Range("J132").Value = 3630.1667
$ws.Range("K132").Value = 8090.667
$ws.Range("L132").Value = 10890.5001
$ws.Range("M132").Value = -5560.667
$ws.Range("N132").Value = -15950.5001

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2386.7
$ws.Range("I82").Value = 2465.65
$ws.Range("J82").Value = 2307.75
$ws.Range("K82").Value = 2465.65
$ws.Range("L82").Value = 2307.75
$ws.Range("M82").Value = -2104.65
$ws.Range("N82").Value = -3029.75

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2386.7
$ws.Range("I85").Value = 2465.65
$ws.Range("J85").Value = 2307.75
$ws.Range("K85").Value = 2465.65
$ws.Range("L85").Value = 2307.75
$ws.Range("M85").Value = -1217.65
$ws.Range("N85").Value = -4803.75

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 228.73334
$ws.Range("I113").Value = 222.91667
$ws.Range("J113").Value = 252
$ws.Range("K113").Value = 668.75001
$ws.Range("L113").Value = 756
$ws.Range("M113").Value = 1501.24999
$ws.Range("N113").Value = -5096

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1628.6487
$ws.Range("I132").Value = 875.9074000000001
$ws.Range("J132").Value = 3661.05
$ws.Range("K132").Value = 2627.7222
$ws.Range("L132").Value = 10983.15
$ws.Range("M132").Value = -97.72220000000016
$ws.Range("N132").Value = -16043.15

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1613.7794
$ws.Range("I136").Value = 860.2857
$ws.Range("J136").Value = 2830.9614
$ws.Range("K136").Value = 2580.8571
$ws.Range("L136").Value = 8492.8842
$ws.Range("M136").Value = -30.85710000000017
$ws.Range("N136").Value = -13592.8842
